# Updated cryptos list on Sat Aug 12 16:40:44 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns for the
# existing coins, and rotates in a new coin (RocketPoolETH) at rank #42 -
# which pushes Quant/Aave/BabyDogeCoin/Aptos/TheSandbox/EnergySwap/RenderToken
# down one spot each (Algorand, formerly the last row, drops off the list).
#
# Every cell below holds a plain-text value (prices use "." as a thousands
# separator in this feed, e.g. "29.479.41", and some look like plain decimals,
# e.g. "1.000" or "6.185") - so each write temporarily forces a Text number
# format to stop Excel's auto-conversion from turning it into a number, then
# restores the cell to the workbook's default "Normal" style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '29.479.41'
Set-TextValue 'E2' '  +0.39%  '
Set-TextValue 'D3' '1.852.67'
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '240.93'
Set-TextValue 'E5' '  +0.81%  '
Set-TextValue 'D6' '0.6306'
Set-TextValue 'E6' '  +0.22%  '
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '0.07685'
Set-TextValue 'E8' '  +1.78%  '
Set-TextValue 'D9' '0.2937'
Set-TextValue 'E9' '  -0.25%  '
Set-TextValue 'D10' '24.74'
Set-TextValue 'E10' '  +0.74%  '
Set-TextValue 'D11' '0.07753'
Set-TextValue 'E11' '  +0.79%  '
Set-TextValue 'D12' '1.863.04'
Set-TextValue 'E12' '  +1.17%  '
Set-TextValue 'D14' '0.6810'
Set-TextValue 'E14' '  +0.41%  '
Set-TextValue 'D15' '0.00001070'
Set-TextValue 'E15' '  +4.56%  '
Set-TextValue 'D16' '83.74'
Set-TextValue 'E16' '  +0.83%  '
Set-TextValue 'D17' '2.122.74'
Set-TextValue 'E17' '  +1.55%  '
Set-TextValue 'D18' '6.185'
Set-TextValue 'E18' '  +0.92%  '
Set-TextValue 'D19' '29.501.51'
Set-TextValue 'E19' '  +0.33%  '
Set-TextValue 'D20' '229.28'
Set-TextValue 'E20' '  +0.48%  '
Set-TextValue 'D21' '12.47'
Set-TextValue 'E21' '  +0.51%  '
Set-TextValue 'E22' '  +0.04%  '
Set-TextValue 'D23' '7.468'
Set-TextValue 'E23' '  +0.23%  '
Set-TextValue 'D25' '157.29'
Set-TextValue 'E25' '  +0.37%  '
Set-TextValue 'E26' '  -0.64%  '
Set-TextValue 'D27' '8.415'
Set-TextValue 'E27' '  +0.79%  '
Set-TextValue 'D28' '17.73'
Set-TextValue 'E28' '  +0.64%  '
Set-TextValue 'D29' '1.334'
Set-TextValue 'E29' '  +5.52%  '
Set-TextValue 'D30' '1.468'
Set-TextValue 'E30' '  +0.62%  '
Set-TextValue 'D31' '0.05692'
Set-TextValue 'E31' '  +1.03%  '
Set-TextValue 'D32' '4.137'
Set-TextValue 'E32' '  +0.44%  '
Set-TextValue 'D33' '4.047'
Set-TextValue 'E33' '  +0.30%  '
Set-TextValue 'D34' '1.856'
Set-TextValue 'E34' '  +1.22%  '
Set-TextValue 'D35' '1.167'
Set-TextValue 'E35' '  +1.06%  '
Set-TextValue 'D36' '0.7086'
Set-TextValue 'E36' '  -0.38%  '
Set-TextValue 'D37' '2.589'
Set-TextValue 'E37' '  -0.15%  '
Set-TextValue 'D38' '2.785'
Set-TextValue 'E38' '  +0.50%  '
Set-TextValue 'D39' '0.01795'
Set-TextValue 'E39' '  -0.72%  '
Set-TextValue 'D40' '1.221.00'
Set-TextValue 'E40' '  -1.51%  '
Set-TextValue 'D41' '6.561'
Set-TextValue 'E41' '  +5.36%  '
Set-TextValue 'D42' '0.9083'
Set-TextValue 'E42' '  +0.81%  '
Set-TextValue 'D43' '1.002'
Set-TextValue 'E43' '  +0.15%  '
Set-TextValue 'B44' 'RocketPoolETH'
Set-TextValue 'C44' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D44' '2.030.69'
Set-TextValue 'E44' '  +1.54%  '
Set-TextValue 'B45' 'Quant'
Set-TextValue 'C45' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D45' '101.75'
Set-TextValue 'E45' '  -0.09%  '
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '66.59'
Set-TextValue 'E46' '  +1.09%  '
Set-TextValue 'B47' 'BabyDogeCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D47' '0.00000000121'
Set-TextValue 'E47' '  +3.10%  '
Set-TextValue 'B48' 'Aptos'
Set-TextValue 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D48' '7.133'
Set-TextValue 'E48' '  +0.43%  '
Set-TextValue 'B49' 'TheSandbox'
Set-TextValue 'C49' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D49' '0.4027'
Set-TextValue 'E49' '  +0.84%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '9.021'
Set-TextValue 'E50' '  +0.76%  '
Set-TextValue 'B51' 'RenderToken'
Set-TextValue 'C51' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D51' '1.690'
Set-TextValue 'E51' '  +0.43%  '
